$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet 1"

# Full replacement data set: header row + 17 product rows
$data = @(
    @("Item name", "Price", "Shop"),
    @("GeForce RTX 3060 Ventus 2X OC 12GB GDDR6 192bit MSI", 290, "RD Electronics"),
    @("GeForce RTX 3060 Gaming OC 8GB GDDR6 128bit rev. 2.0 Gigabyte", 292, "RD Electronics"),
    @("GeForce RTX 3060 Gaming OC 8GB GDDR6 128bit Gigabyte", 292.63, "RD Electronics"),
    @("GeForce RTX 3060 WindForce OC 12GB GDDR6 192bit Gigabyte", 300.18, "RD Electronics"),
    @("GeForce RTX 3060 Gaming OC 12GB GDDR6 192bit rev. 2.0 LHR Gigabyte", 302.99, "RD Electronics"),
    @("GeForce RTX 3060 Dual OC V2 LHR 12GB GDDR6 192bit Asus", 317.9, "RD Electronics"),
    @("GeForce RTX Extreme Gaming 3060 12GB GDDR6 192bit Biostar", 317.9, "RD Electronics"),
    @("GeForce RTX 3060 V2 OC 12GB GDDR6 192bit Asus", 329, "RD Electronics"),
    @("GeForce RTX 3060 Twin Twin Edge 12GB GDDR6 192bit ZOTAC", 334.98, "RD Electronics"),
    @("GeForce RTX 3060 Twin X2 OC LHR 12GB GDDR6 192bit Inno3D", 341, "RD Electronics"),
    @("GeForce RTX 3060 Twin X2 12GB GDDR6 192bit Inno3D", 351, "RD Electronics"),
    @("GeForce RTX 3060 Ti Extreme Gaming 8GB GDDR6X 256bit Biostar", 361, "RD Electronics"),
    @("Ryzen 5 5500 16GB 250SSD RTX3060 12GB no-OS INTOP", 699.99, "RD Electronics"),
    @("Ryzen 5 5500 16GB 500SSD RTX3060 12GB no-OS INTOP", 711.99, "RD Electronics"),
    @("Ryzen 5 5500 32GB 250SSD RTX3060 12GB no-OS INTOP", 728.99, "RD Electronics"),
    @("Ryzen 5 5600X 16GB 250SSD RTX3060 12GB no-OS INTOP", 734.99, "RD Electronics"),
    @("Blade RZ09-0421EED3-R3E1 15.6 240hz i7-12800H 16GB 1SSD RTX3060 EN W11 Razer", 3340, "RD Electronics")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
